$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Parameter_Input"
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Parameter_Input")

# Row 2 - existing series config, ticker swapped to ODSACBW027SBOG (Other Deposits)
$ws1.Range("B2").Value = "ODSACBW027SBOG"
$ws1.Range("D2").Value = "Year on year % change"
$ws1.Range("E2").Value = "black"
$ws1.Range("F2").Value = "odl"
$ws1.Range("H2").Value = "Other Deposits, All Commercial Banks"
$ws1.Range("I2").Value = "Trillions of U.S `$"
$ws1.Range("N2").Value = $null

# Row 3 - existing series config, switched to BTCUSD,INDEX (Bitcoin)
$ws1.Range("B3").Value = "BTCUSD,INDEX"
$ws1.Range("D3").Value = "Year on year % change"
$ws1.Range("E3").Value = "orangered"
$ws1.Range("F3").Value = "btc"
$ws1.Range("H3").Value = "Bitcoin"
$ws1.Range("I3").Value = $null
$ws1.Range("N3").Value = 2.5

# Row 4 - newly populated with S&P 500 series
$ws1.Range("B4").Value = "^GSPC"
$ws1.Range("C4").Value = "load"
$ws1.Range("D4").Value = "Year on year % change"
$ws1.Range("E4").Value = "blue"
$ws1.Range("F4").Value = "SPX"
$ws1.Range("G4").Value = "log"
$ws1.Range("H4").Value = "S & P 500"
$ws1.Range("I4").Value = "USD"

# Row 5 - newly populated with Global M2 (top 50) series
$ws1.Range("B5").Value = "Top50GM2"
$ws1.Range("C5").Value = "load"
$ws1.Range("D5").Value = "Year on year % change"
$ws1.Range("E5").Value = "green"
$ws1.Range("F5").Value = "GM2"
$ws1.Range("G5").Value = "log"
$ws1.Range("H5").Value = "Global M2 (top 50)"
$ws1.Range("N5").Value = 2.5

# Row 8 - StartDate
$ws1.Range("B8").Value = "2010-01-01"

# Row 12 - CHART TITLE
$ws1.Range("B12").Value = "Monetary Aggregates, Bitcoin & Equities."

# ------------------------------------------------------------------
# Sheet "TickerList"
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("TickerList")

# Row 52 - M2SL catalog entry
$ws2.Range("B52").Value = "M2SL"
$ws2.Range("C52").Value = "load"
$ws2.Range("D52").Value = "Unaltered"
$ws2.Range("E52").Value = "black"
$ws2.Range("F52").Value = "M2"
$ws2.Range("G52").Value = "log"
$ws2.Range("H52").Value = "M2 money supply (U.S)"
$ws2.Range("I52").Value = "Trillions of U.S `$"
$ws2.Range("J52").Value = 1000

# Row 53 - ODSACBW027SBOG catalog entry
$ws2.Range("B53").Value = "ODSACBW027SBOG"
$ws2.Range("C53").Value = "load"
$ws2.Range("D53").Value = "Unaltered"
$ws2.Range("E53").Value = "green"
$ws2.Range("F53").Value = "odl"
$ws2.Range("G53").Value = "log"
$ws2.Range("H53").Value = "Other Deposits, All Commercial Banks"
$ws2.Range("I53").Value = "Trillions of U.S `$"
$ws2.Range("J53").Value = 1000

# Row 54 - MABMM301USM189S catalog entry (includes a formula for J54)
$ws2.Range("B54").Value = "MABMM301USM189S"
$ws2.Range("C54").Value = "load"
$ws2.Range("D54").Value = "Unaltered"
$ws2.Range("E54").Value = "red"
$ws2.Range("F54").Value = "M3"
$ws2.Range("G54").Value = "log"
$ws2.Range("H54").Value = "M3 - broad money "
$ws2.Range("I54").Value = "Trillions of U.S `$"
$ws2.Range("J54").Formula = "=10^12"

# Row 55 - ^GSPC catalog entry
$ws2.Range("B55").Value = "^GSPC"
$ws2.Range("C55").Value = "load"
$ws2.Range("D55").Value = "Unaltered"
$ws2.Range("E55").Value = "blue"
$ws2.Range("F55").Value = "SPX"
$ws2.Range("G55").Value = "log"
$ws2.Range("H55").Value = "S & P 500"
$ws2.Range("I55").Value = "USD"
